# Aguinaldo.xlsx — "Add files via upload"
# Fill in the Performance Task / Written Work grade rows for the two
# students plus the "Highest Possible Grade" reference row, and touch up
# the sheet view (zoom + selection) to match the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aguinaldo")

# --- Row 3: PONS, Luke Alexander R -----------------------------------
$ws.Range("A3").Value = "PONS, Luke Alexander R"
$ws.Range("B3:U3").Value = 10

# --- Row 4: GIDA, Ericson Virgile F -----------------------------------
$ws.Range("A4").Value = "GIDA, Ericson Virgile F"
$ws.Range("B4:U4").Value = 10

# --- Row 5: Highest Possible Grade ------------------------------------
$ws.Range("A5").Value = "Highest Possible Grade"
$ws.Range("B5:U5").Value = 10

# --- View state: zoom to 90% and leave the cursor on O12 --------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 90
$ws.Range("O12").Select()
